$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# New header cell CC1 - copy formatting (style) from CB1, then set the text value
$ws.Range("CB1").Copy()
$ws.Range("CC1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("CC1").Value = "02-sep"

# New data values for CC2:CC25 (plain numbers, no special style)
$values = @(33.61, 36.79, 52.57, 40.16, 42.55, 12.47, 33.44, 34.28, 34.8, 39, 7.96, 1.2, 1.5, 0.65, 0.21, 1.5, 8.210000000000001, 14, 32.89, 84, 53.46, 34.98, 42.89, 31.06)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 81).Value = $values[$i]
}

$wb.Save()
